$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12, shifting existing rows 12-16 down to 13-17
$ws.Rows.Item(12).Insert()

# Populate the new row 12 with data (columns A, B, C, E, F, G, Q, R carry over
# the same values as the old row 12 had; D, H, I, J, K, L, M, N, O, P are new)
$ws.Range("A12").Value = 12
$ws.Range("B12").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C12").Value = "Metropolitana"
$ws.Range("D12").Value = 44488
$ws.Range("E12").Value = 13
$ws.Range("F12").Value = 100112028
$ws.Range("G12").Value = "Sandia"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 150
$ws.Range("K12").Value = 800
$ws.Range("L12").Value = 800
$ws.Range("M12").Value = 800
$ws.Range("N12").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O12").Value = "Perú"
$ws.Range("P12").Value = 800
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = "Hortaliza"
